# Auto-generated: apply cryptos-list price/volume refresh from scrape run
$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

$ws.Range("D2").Value = "66.954.57"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "3.449.83"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "579.84"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").Value = "188.89"
$ws.Range("E6").Value = "  +8.62%  "
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "3.440.62"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").Value = "0.645"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Formula = "'56.90"
$ws.Range("E12").Value = "  +6.23%  "
$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "9.45"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "3.989.25"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Formula = "'18.80"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").Value = "3.440.95"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").Value = "66.848.84"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Formula = "'12.10"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "1.03"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("D22").Value = "483.31"
$ws.Range("E22").Value = "  +6.31%  "
$ws.Range("D23").Value = "5.33"
$ws.Range("E23").Value = "  +8.25%  "
$ws.Range("E24").Value = "  +24.08%  "
$ws.Range("D25").Value = "4.34"
$ws.Range("E25").Value = "  +7.06%  "
$ws.Range("D26").Value = "89.45"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "2.98"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("D29").Value = "9.06"
$ws.Range("E29").Value = "  +4.71%  "
$ws.Range("D30").Value = "31.27"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Formula = "'7.40"
$ws.Range("E31").Value = "  +13.37%  "
$ws.Range("D32").Value = "11.79"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "598.97"
$ws.Range("E33").Value = "  +4.11%  "
$ws.Range("E35").Value = "  +4.34%  "
$ws.Range("D36").Value = "0.149"
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("D37").Formula = "'1.00"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "36.85"
$ws.Range("E38").Value = "  +3.44%  "
$ws.Range("D39").Value = "0.389"
$ws.Range("E39").Value = "  +5.00%  "
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("D41").Value = "0.0₃0754"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").Value = "3.236.07"
$ws.Range("E42").Value = "  +5.35%  "
$ws.Range("D43").Formula = "'2.90"
$ws.Range("E43").Value = "  +5.33%  "
$ws.Range("D44").Value = "0.0431"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").Value = "2.84"
$ws.Range("E45").Value = "  +26.33%  "
$ws.Range("D46").Value = "2.54"
$ws.Range("E46").Value = "  +4.09%  "
$ws.Range("D47").Value = "3.24"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.135"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").Value = "3.36"
$ws.Range("E49").Value = "  +16.58%  "
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "8.65"
$ws.Range("E51").Value = "  +4.60%  "
